$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Build-RunsXml($texts, $rPr) {
    # Builds a concatenation of <w:r> elements, one per text chunk, all
    # sharing the same <w:rPr> block. xml:space="preserve" is only added
    # when the chunk has leading/trailing whitespace (matches Word's own
    # habit of omitting it when not needed).
    $sb = New-Object System.Text.StringBuilder
    foreach ($t in $texts) {
        [void]$sb.Append("<w:r>")
        [void]$sb.Append($rPr)
        $needsPreserve = ($t.Length -eq 0) -or ($t.Substring(0,1) -eq ' ') -or ($t.Substring($t.Length-1,1) -eq ' ')
        if ($needsPreserve) {
            [void]$sb.Append('<w:t xml:space="preserve">')
        } else {
            [void]$sb.Append('<w:t>')
        }
        [void]$sb.Append($t)
        [void]$sb.Append('</w:t></w:r>')
    }
    return $sb.ToString()
}

function Wrap-Package($bodyInnerXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$rPrNormal = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

# ---------------------------------------------------------------------------
# Change 1: expand the "ROC AUC score of approximately 0.63 ... training
# set." sentence into the longer, multi-sentence explanation.
# ---------------------------------------------------------------------------

$target1 = "This model achieves the ROC AUC score of approximately 0.63 on the validation test of 20% the training set."

$found1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq $target1) {
        $found1 = $p
        break
    }
}
if ($found1 -eq $null) { throw "Change 1 target paragraph not found" }

$chunks1 = @(
    "This model achieves the ROC AUC score of approximately 0.63",
    "-0.65",
    " ",
    "on the validation test of 20% the training set",
    " when the prediction label is converted to 0-1 labels with 0.5 threshold",
    ".",
    " When calculating ",
    "ROC AUC score",
    " with neural network ",
    "probability ",
    "output (not converted to 0-1 labels), ",
    "ROC AUC",
    " score is ",
    "around 0.70-",
    " 0.71",
    " (depending on specific train-test split)",
    "."
)

$runsXml1 = Build-RunsXml $chunks1 $rPrNormal
$pkg1 = Wrap-Package ("<w:p>" + $runsXml1 + "</w:p>")

$pRange1 = $found1.Range
$subRange1 = $d.Range($pRange1.Start, $pRange1.End - 1)
$subRange1.InsertXML($pkg1)

# ---------------------------------------------------------------------------
# Change 2: "Obviously a ROC AUC score of 0.63 is not ideal..." becomes
# "A ROC AUC score of around 0.7 is not ideal...", with the bookmarked
# _GoBack location (and the lastRenderedPageBreak marker) moving to sit
# right after the new opening word "A".
# ---------------------------------------------------------------------------

$target2 = "Obviously a ROC AUC score of 0.63 is not ideal. I suspect that with some parameter tuning, I can increase the score by a few percentages. However, for a significantly better score, this problem will probably require a new approach. There are a few methods I have read about but have not tried due to time constraint. I will keep exploring this problem, and I would like to receive any suggestion for a solution."

$found2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq $target2) {
        $found2 = $p
        break
    }
}
if ($found2 -eq $null) { throw "Change 2 target paragraph not found" }

# Remove the pre-existing _GoBack bookmark so re-adding it below doesn't
# leave a stray duplicate behind.
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rPrWithBreak = $rPrNormal + '<w:lastRenderedPageBreak/>'
$run2First = "<w:r>" + $rPrWithBreak + "<w:t>A</w:t></w:r>"

$chunks2Rest = @(
    " ",
    "ROC AUC score of ",
    "around 0.7",
    " is not ideal. I suspect that with some parameter tuning, I can increase the score by a few percentages. ",
    "However, for a significantly better score, this problem will probably require a new approach. There are a few methods I have read about but have not tried due to time constraint. I will keep exploring this problem, and I would like to receive any suggestion for a solution."
)
$runsXml2Rest = Build-RunsXml $chunks2Rest $rPrNormal

$pkg2 = Wrap-Package ("<w:p>" + $run2First + $runsXml2Rest + "</w:p>")

$pRange2 = $found2.Range
$subRange2 = $d.Range($pRange2.Start, $pRange2.End - 1)
$subRange2.InsertXML($pkg2)

# Re-insert the _GoBack bookmark right after the new leading "A".
$pRange2b = $found2.Range
$bmPos = $pRange2b.Start + 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "Done"
